# Fix spaces in dataset filenames
# The "datasets" worksheet (sheet3) has a column H containing filenames of the
# form "MedBFM4_<var>_6869436a-80f4-4c6d-954b-a730b348d7ce .tif" (note the
# stray space before the extension). Remove that stray space for rows 29-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("datasets")

# Mirrors the order Excel's "Replace All" walked the cells (wrapping from the
# previously active cell H45 back around to H29) so the rebuilt shared-string
# table lands in the same order as the author's saved file.
$rows = @(45, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 8)  # column H
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace ' \.tif$', '.tif'
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# Reflect the author's final view state: active cell I14 and top-left cell F1.
$ws.Activate()
$ws.Range("I14").Select()
$excel.ActiveWindow.ScrollColumn = 6
